$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.941.48"
$ws.Range("E2").Value = "'  -0.48%  "
$ws.Range("D3").Value = "'1.703.34"
$ws.Range("E3").Value = "'  -0.73%  "
$ws.Range("E4").Value = "'  -0.56%  "
$ws.Range("D5").Value = "'316.59"
$ws.Range("E5").Value = "'  -0.59%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "'  -0.35%  "
$ws.Range("D7").Value = "'0.4032"
$ws.Range("E7").Value = "'  +1.49%  "
$ws.Range("D8").Value = "'0.4061"
$ws.Range("E8").Value = "'  -1.74%  "
$ws.Range("D9").Value = "'1.003"
$ws.Range("E9").Value = "'  -0.56%  "
$ws.Range("E10").Value = "'  -3.90%  "
$ws.Range("D11").Value = "'53.72"
$ws.Range("E11").Value = "'  +1.42%  "
$ws.Range("D12").Value = "'0.08821"
$ws.Range("E12").Value = "'  -1.25%  "
$ws.Range("D13").Value = "'26.01"
$ws.Range("E13").Value = "'  +3.93%  "
$ws.Range("D14").Value = "'7.478"
$ws.Range("E14").Value = "'  -3.12%  "
$ws.Range("D15").Value = "'8.039"
$ws.Range("E15").Value = "'  -1.42%  "
$ws.Range("D16").Value = "'0.00001351"
$ws.Range("E16").Value = "'  -3.18%  "
$ws.Range("D17").Value = "'1.728.19"
$ws.Range("E17").Value = "'  +1.08%  "
$ws.Range("D18").Value = "'96.35"
$ws.Range("E18").Value = "'  -4.48%  "
$ws.Range("D19").Value = "'0.07158"
$ws.Range("E19").Value = "'  +0.19%  "
$ws.Range("D20").Value = "'20.95"
$ws.Range("E20").Value = "'  +3.33%  "
$ws.Range("D21").Value = "'7.246"
$ws.Range("E21").Value = "'  -3.14%  "
$ws.Range("D22").Value = "'1.004"
$ws.Range("E22").Value = "'  -0.58%  "
$ws.Range("D23").Value = "'14.51"
$ws.Range("E23").Value = "'  -0.10%  "
$ws.Range("D24").Value = "'24.917.58"
$ws.Range("E24").Value = "'  -0.60%  "
$ws.Range("D25").Value = "'2.332"
$ws.Range("E25").Value = "'  -0.78%  "
$ws.Range("D26").Value = "'2.892"
$ws.Range("E26").Value = "'  -7.02%  "
$ws.Range("D27").Value = "'6.503"
$ws.Range("E27").Value = "'  +24.64%  "
$ws.Range("D28").Value = "'23.06"
$ws.Range("E28").Value = "'  -0.14%  "
$ws.Range("D29").Value = "'165.78"
$ws.Range("E29").Value = "'  +0.08%  "
$ws.Range("D30").Value = "'145.29"
$ws.Range("E30").Value = "'  +4.18%  "
$ws.Range("D31").Value = "'8.234"
$ws.Range("E31").Value = "'  -6.62%  "
$ws.Range("D32").Value = "'1.916.04"
$ws.Range("E32").Value = "'  +0.76%  "
$ws.Range("D33").Value = "'2.232"
$ws.Range("E33").Value = "'  +13.00%  "
$ws.Range("D34").Value = "'0.08864"
$ws.Range("E34").Value = "'  -1.59%  "
$ws.Range("D35").Value = "'7.412"
$ws.Range("E35").Value = "'  -5.12%  "
$ws.Range("D36").Value = "'0.03203"
$ws.Range("E36").Value = "'  +6.88%  "
$ws.Range("D37").Value = "'1.016"
$ws.Range("E37").Value = "'  -5.97%  "
$ws.Range("D38").Value = "'0.2842"
$ws.Range("E38").Value = "'  +1.63%  "
$ws.Range("B39").Value = "'FraxShare"
$ws.Range("C39").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'10.85"
$ws.Range("E39").Value = "'  -1.93%  "
$ws.Range("B40").Value = "'TheSandbox"
$ws.Range("C40").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.8385"
$ws.Range("E40").Value = "'  +2.25%  "
$ws.Range("D41").Value = "'0.09356"
$ws.Range("E41").Value = "'  +0.98%  "
$ws.Range("D42").Value = "'14.07"
$ws.Range("E42").Value = "'  -4.00%  "
$ws.Range("D43").Value = "'1.478"
$ws.Range("E43").Value = "'  -0.40%  "
$ws.Range("D44").Value = "'17.55"
$ws.Range("E44").Value = "'  +6.19%  "
$ws.Range("D45").Value = "'2.717"
$ws.Range("E45").Value = "'  +2.70%  "
$ws.Range("D46").Value = "'0.7433"
$ws.Range("E46").Value = "'  +0.28%  "
$ws.Range("D47").Value = "'4.248"
$ws.Range("E47").Value = "'  -0.97%  "
$ws.Range("D48").Value = "'1.389"
$ws.Range("E48").Value = "'  +2.78%  "
$ws.Range("E49").Value = "'  -0.38%  "
$ws.Range("D50").Value = "'142.24"
$ws.Range("E50").Value = "'  +1.10%  "
$ws.Range("D51").Value = "'0.08349"
$ws.Range("E51").Value = "'  +2.96%  "
